$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header cells J1:L1 (birthday / about_me / picture), styled like the
#        existing header row (bold, centered, bordered) by pasting A1's format. ---
$ws.Range("J1").Value = "birthday"
$ws.Range("K1").Value = "about_me"
$ws.Range("L1").Value = "picture"
$ws.Range("A1").Copy()
$ws.Range("J1:L1").PasteSpecial(-4122)

# --- 2. Blank J2:L4 (new columns for existing rows 2-4) - set to empty string
#        then paste the (unstyled) format of A2 so the cells persist as blank
#        rather than being dropped entirely. ---
$ws.Range("J2:L4").Value = ""
$ws.Range("A2").Copy()
$ws.Range("J2:L4").PasteSpecial(-4122)

# --- 3. Update row 4 in place: swap Parth Sharadrao's profile for Vedant Dinkar's. ---
$ws.Range("B4").Value = "Vedant Dinkar"
$ws.Range("D4").Value = "Undergraduate Student of CSE at Indian Institute of Technology, Indore"
$ws.Range("E4").Value = "['Google Developer Student Clubs', 'The Debating Society IIT Indore', 'Model United Nations, IIT Indore']"
$ws.Range("F4").Value = "['Volunteer', 'Volunteer', 'Executive Board Affairs Member']"
$ws.Range("G4").Value = "https://www.linkedin.com/in/vedant-dinkar-a6a4301b9/"
$ws.Range("I4").Value = "cse220001078@iiti.ac.in"

# --- 4. Append new row 5, duplicating Vedant Dinkar's profile and adding the new
#        birthday / about_me fields. ---
$ws.Range("A5").Value = "India"
$ws.Range("B5").Value = "Vedant Dinkar"
$ws.Range("D5").Value = "Undergraduate Student of CSE at Indian Institute of Technology, Indore"
$ws.Range("E5").Value = "['Google Developer Student Clubs', 'The Debating Society IIT Indore', 'Model United Nations, IIT Indore']"
$ws.Range("F5").Value = "['Volunteer', 'Volunteer', 'Executive Board Affairs Member']"
$ws.Range("G5").Value = "https://www.linkedin.com/in/vedant-dinkar-a6a4301b9/"
$ws.Range("H5").Value = "Computer Science"
$ws.Range("I5").Value = "cse220001078@iiti.ac.in"

# J5 ("2004-08-10") looks like a date, and a plain .Value assignment would get
# auto-converted to a date serial number by Excel's type inference. Route it
# through a formula-literal + paste-values round trip (via a scratch cell well
# outside the sheet's used range) so it lands as plain text instead, matching
# the source data. The scratch cell is cleared afterwards so it leaves no trace.
$ws.Range("Z1").Formula = '="2004-08-10"'
$ws.Range("Z1").Copy()
$ws.Range("J5").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

$ws.Range("K5").Value = "The power is within us..."

# C5 (last_name, unused) and L5 (new picture column) stay blank - same
# blank-cell trick as used for columns J:L in rows 2-4 above, so the cells
# persist instead of being dropped.
$ws.Range("C5").Value = ""
$ws.Range("L5").Value = ""
$ws.Range("A2").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("L5").PasteSpecial(-4122)
